$d = $word.ActiveDocument

# The document's header/footer logo pictures had their display "name"
# (wp:docPr / pic:cNvPr name="...") swapped between the two Pearson
# logos (image1.png <-> image2.png) and renumbered for the BTec logo
# (image2.jpg -> image1.jpg). The underlying embedded picture bytes
# and everything else (ids, alt-text/descr, size, position) stay the
# same - only the "name" metadata changes.
#
# InlineShape has no settable Name in the Word object model, so each
# picture is briefly converted to a floating Shape (where .Name is
# writable), renamed, then converted back to an inline picture.

foreach ($sec in $d.Sections) {

    # Headers: wdHeaderFooterPrimary=1, wdHeaderFooterFirstPage=2, wdHeaderFooterEvenPage=3
    for ($h = 1; $h -le $sec.Headers.Count; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
                $ils = $hdr.Range.InlineShapes.Item($i)
                $oldName = $ils.AlternativeText
                $shp = $ils.ConvertToShape()
                if ($oldName -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
                $shp.ConvertToInlineShape()
            }
        }
    }

    # Footers: same index scheme as Headers
    for ($f = 1; $f -le $sec.Footers.Count; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
                $ils = $ftr.Range.InlineShapes.Item($i)
                $oldName = $ils.AlternativeText
                $shp = $ils.ConvertToShape()
                if ($oldName -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
                $shp.ConvertToInlineShape()
            }
        }
    }
}

Write-Host "Renamed header/footer logo pictures"
